# Insert a new data row at row 449 (pushing existing rows 449:538 down to 450:539)
# and populate it with the new price-record values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 449, shifting everything below it down.
$ws.Rows.Item(449).Insert(-4121, 0)   # -4121 = xlShiftDown, 0 = xlFormatFromLeftOrAbove

# Fill in the new row's values (columns A-T)
$ws.Cells.Item(449, 1).Value = 10
$ws.Cells.Item(449, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(449, 3).Value = "La Araucanía"
$ws.Cells.Item(449, 4).Value = 44511
$ws.Cells.Item(449, 5).Value = 9
$ws.Cells.Item(449, 6).Value = "Fruta"
$ws.Cells.Item(449, 7).Value = 100102
$ws.Cells.Item(449, 8).Value = "Cítricos"
$ws.Cells.Item(449, 9).Value = 100102005
$ws.Cells.Item(449, 10).Value = "Naranja"
$ws.Cells.Item(449, 11).Value = "Navel Late"
$ws.Cells.Item(449, 12).Value = "Primera"
$ws.Cells.Item(449, 13).Value = 220
$ws.Cells.Item(449, 14).Value = 9000
$ws.Cells.Item(449, 15).Value = 10000
$ws.Cells.Item(449, 16).Value = 9705
$ws.Cells.Item(449, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(449, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(449, 19).Value = 647
$ws.Cells.Item(449, 20).Value = 15

# Match the date-cell style used by the rest of column D
$ws.Cells.Item(449, 4).NumberFormat = $ws.Cells.Item(450, 4).NumberFormat
